$d = $word.ActiveDocument

function Replace-All($find, $replace) {
    [void]$d.Content.Find.Execute($find, $false, $false, $false, $false, $false,
                                   $true, 1, $false, $replace, 2)
}

# Title
Replace-All "ARTIFICIAL INTELLIGENCE AND MACHINE LEARNING" "PRODUCT DEVELOPMENT"

# More specific, longer phrases first to avoid clobbering shorter substrings
Replace-All "Project Name: Product and Machine Learning Implementation Initiative" "Project Name: Product Development and Product Innovation Implementation Initiative"
Replace-All "Industry: Product and Machine Learning" "Industry: Product Development and Product Innovation"
Replace-All "Industry Focus: Product and Machine Learning" "Industry Focus: Product Development and Product Innovation"
Replace-All "initiative for Product and Machine Learning to achieve" "initiative for Product Development and Product Innovation to achieve"

Replace-All "Project Type: Product Implementation" "Project Type: Product Development Implementation"
Replace-All "strategic Product Implementation initiative" "strategic Product Development Implementation initiative"
Replace-All "through Product Implementation capabilities." "through Product Development Implementation capabilities."
Replace-All "comprehensive Product Implementation solution" "comprehensive Product Development Implementation solution"

Replace-All "ML Engineers" "Product Engineers"
Replace-All "Compliance Officers" "Quality Assurance Managers"

Replace-All "MLflow" "Productflow"
Replace-All "Cloud ML platforms" "Cloud Product platforms"

Replace-All "Financial Justification:" "Product Justification:"
Replace-All "Financial: Budget overruns, cost escalation, ROI delays" "Product: Budget overruns, cost escalation, ROI delays"
Replace-All "ML Platform Licensing: `$840,000" "Product Platform Licensing: `$840,000"
Replace-All "Financial: Break-even within 30 months, 250%+ ROI within 3 years" "Product: Break-even within 30 months, 250%+ ROI within 3 years"
Replace-All "Financial review and budget allocation approval" "Product review and budget allocation approval"

# Remove the page-break paragraph: it currently holds a single run
# containing only a page-type <w:br/>. The target document keeps the
# paragraph but with a single empty run and no break. Locate it precisely
# by searching for the page-break special character (^m) within each
# paragraph's own range.
$pageBreakXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
foreach ($p in $d.Paragraphs) {
    $rng = $p.Range.Duplicate
    $hasBreak = $rng.Find.Execute("^m", $false, $false, $false, $false, $false,
                                   $true, 1, $false, "", 0)
    if ($hasBreak) {
        [void]$p.Range.InsertXML($pageBreakXml)
        break
    }
}


